# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# currentAveragePrice / LevePrice / LeveProfit figures from the scheduled
# market-data runner. A few rows also gain or lose their HQ/NQ profit cell
# (LeveProfitNQ / LeveProfitHQ) entirely when the corresponding HQ/NQ price
# crosses to/from zero, matching how the source generator omits zero cells.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 343.6154
$ws.Range("J33").Value = 82.5
$ws.Range("L33").Value = 82.5
$ws.Range("N33").Value = -540.5
$ws.Range("H74").Value = 4208.5
$ws.Range("I74").Value = 4208.5
$ws.Range("K74").Value = 4208.5
$ws.Range("M74").Value = -3272.5
$ws.Range("H77").Value = 4208.5
$ws.Range("I77").Value = 4208.5
$ws.Range("K77").Value = 21042.5
$ws.Range("M77").Value = -16362.5
$ws.Range("H98").Value = 1168.0834
$ws.Range("I98").Value = 1047
$ws.Range("K98").Value = 1047
$ws.Range("M98").Value = 451
$ws.Range("H122").Value = 1168.0834
$ws.Range("I122").Value = 1047
$ws.Range("K122").Value = 3141
$ws.Range("M122").Value = -691
$ws.Range("H132").Value = 6555.1333
$ws.Range("I132").Value = 6766.5
$ws.Range("K132").Value = 20299.5
$ws.Range("M132").Value = -17769.5
$ws.Range("H135").Value = 825.6667
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1232.421
$ws.Range("I137").Value = 1299.8125
$ws.Range("K137").Value = 3899.4375
$ws.Range("M137").Value = -1349.4375

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3496
$ws.Range("I32").Value = 2154.138
$ws.Range("K32").Value = 2154.138
$ws.Range("M32").Value = -1867.138
$ws.Range("H61").Value = 1575.3
$ws.Range("I61").Value = 1528.1111
$ws.Range("K61").Value = 1528.1111
$ws.Range("M61").Value = -1316.1111
$ws.Range("H74").Value = 702.9167
$ws.Range("I74").Value = 630.4545000000001
$ws.Range("K74").Value = 630.4545000000001
$ws.Range("M74").Value = 243.5454999999999
$ws.Range("H77").Value = 702.9167
$ws.Range("I77").Value = 630.4545000000001
$ws.Range("K77").Value = 3152.2725
$ws.Range("M77").Value = 1215.7275
$ws.Range("H122").Value = 1686.9
$ws.Range("I122").Value = 1712.5264
$ws.Range("K122").Value = 5137.5792
$ws.Range("M122").Value = -2687.5792
$ws.Range("H136").Value = 1575.3
$ws.Range("I136").Value = 1528.1111
$ws.Range("K136").Value = 4584.3333
$ws.Range("M136").Value = -2034.3333

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 426.53845
$ws.Range("I94").Value = 332.8
$ws.Range("K94").Value = 332.8
$ws.Range("M94").Value = 118.2
$ws.Range("H105").Value = 3847.0417
$ws.Range("I105").Value = 3029.2856
$ws.Range("J105").Value = 4991.9
$ws.Range("K105").Value = 3029.2856
$ws.Range("L105").Value = 4991.9
$ws.Range("M105").Value = -1282.2856
$ws.Range("N105").Value = -8485.9
$ws.Range("H134").Value = 3780.6365
$ws.Range("I134").Value = 3908.7
$ws.Range("K134").Value = 11726.1
$ws.Range("M134").Value = -9191.099999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 29959
$ws.Range("J20").Value = 29959
$ws.Range("L20").Value = 29959
$ws.Range("N20").Value = -30431
$ws.Range("H22").Value = 25447
$ws.Range("I22").Value = 595.6667
$ws.Range("K22").Value = 595.6667
$ws.Range("M22").Value = -245.6667
$ws.Range("H30").Value = 29959
$ws.Range("J30").Value = 29959
$ws.Range("L30").Value = 29959
$ws.Range("N30").Value = -30141
$ws.Range("H58").Value = 5629.7144
$ws.Range("I58").Value = 5734.6665
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 5734.6665
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -5531.6665
$ws.Range("N58").Value = -5406
$ws.Range("H86").Value = 8000
$ws.Range("J86").Value = 8000
$ws.Range("L86").Value = 8000
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 8000
$ws.Range("J89").Value = 8000
$ws.Range("L89").Value = 40000
$ws.Range("N89").Value = -51232
$ws.Range("H128").Value = 29959
$ws.Range("J128").Value = 29959
$ws.Range("L128").Value = 29959
$ws.Range("N128").Value = -39919
$ws.Range("H136").Value = 5629.7144
$ws.Range("I136").Value = 5734.6665
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 17203.9995
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -14653.9995
$ws.Range("N136").Value = -20100

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 252.2
$ws.Range("I13").Value = 37
$ws.Range("J13").Value = 575
$ws.Range("K13").Value = 111
$ws.Range("L13").Value = 1725
$ws.Range("M13").Value = 57
$ws.Range("N13").Value = -2061
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H109").Value = 1800.7273
$ws.Range("I109").Value = 934.6667
$ws.Range("K109").Value = 2804.0001
$ws.Range("M109").Value = -1764.0001
$ws.Range("H113").Value = 1785.8334
$ws.Range("J113").Value = 1693.3334
$ws.Range("L113").Value = 5080.0002
$ws.Range("N113").Value = -9420.0002
$ws.Range("H131").Value = 2750
$ws.Range("I131").Value = 2000
$ws.Range("J131").Value = 3000
$ws.Range("K131").Value = 6000
$ws.Range("L131").Value = 9000
$ws.Range("M131").Value = -960
$ws.Range("N131").Value = -19080

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 181.85715
$ws.Range("I2").Value = 263.5
$ws.Range("J2").Value = 73
$ws.Range("K2").Value = 263.5
$ws.Range("L2").Value = 73
$ws.Range("M2").Value = -150.5
$ws.Range("N2").Value = -299
$ws.Range("H70").Value = 7396.857
$ws.Range("I70").Value = 6355.6
$ws.Range("K70").Value = 6355.6
$ws.Range("M70").Value = -6085.6
$ws.Range("H73").Value = 7396.857
$ws.Range("I73").Value = 6355.6
$ws.Range("K73").Value = 6355.6
$ws.Range("M73").Value = -5419.6
$ws.Range("H80").Value = 2894.75
$ws.Range("I80").Value = 2821
$ws.Range("K80").Value = 2821
$ws.Range("M80").Value = -1823
$ws.Range("H83").Value = 2894.75
$ws.Range("I83").Value = 2821
$ws.Range("K83").Value = 14105
$ws.Range("M83").Value = -9113
$ws.Range("H102").Value = 3121.375
$ws.Range("I102").Value = 3993.5
$ws.Range("K102").Value = 3993.5
$ws.Range("M102").Value = -2371.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5843.6665
$ws.Range("I7").Value = 5843.6665
$ws.Range("K7").Value = 5843.6665
$ws.Range("M7").Value = -5731.6665
$ws.Range("H31").Value = 9681.25
$ws.Range("I31").Value = 5624.5
$ws.Range("J31").Value = 11709.625
$ws.Range("K31").Value = 5624.5
$ws.Range("L31").Value = 11709.625
$ws.Range("M31").Value = -5376.5
$ws.Range("N31").Value = -12205.625
$ws.Range("H40").Value = 1418.7693
$ws.Range("I40").Value = 1418.7693
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1418.7693
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1282.7693
$ws.Range("N40").ClearContents()
$ws.Range("H82").Value = 1359.8
$ws.Range("I82").Value = 1599.6666
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 1599.6666
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -1238.6666
$ws.Range("N82").Value = -1722
$ws.Range("H85").Value = 1359.8
$ws.Range("I85").Value = 1599.6666
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 1599.6666
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = -351.6666
$ws.Range("N85").Value = -3496
$ws.Range("H126").Value = 5843.6665
$ws.Range("I126").Value = 5843.6665
$ws.Range("K126").Value = 17530.9995
$ws.Range("M126").Value = -15060.9995
$ws.Range("H128").Value = 76998.60000000001
$ws.Range("J128").Value = 76998.60000000001
$ws.Range("L128").Value = 76998.60000000001
$ws.Range("N128").Value = -86958.60000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3101.2285
$ws.Range("I122").Value = 3194.7932
$ws.Range("K122").Value = 9584.3796
$ws.Range("M122").Value = -7134.3796
$ws.Range("H130").Value = 53266
$ws.Range("J130").Value = 53266
$ws.Range("L130").Value = 53266
$ws.Range("N130").Value = -63306
$ws.Range("H136").Value = 4615.8057
$ws.Range("I136").Value = 4664.0586
$ws.Range("J136").Value = 4572.6313
$ws.Range("K136").Value = 13992.1758
$ws.Range("L136").Value = 13717.8939
$ws.Range("M136").Value = -11442.1758
$ws.Range("N136").Value = -18817.8939
